$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes the existing May/April/March data
# down by one), so the new "day 3 of June 2025" entry lands right after the
# existing June 1st/2nd rows.
$ws.Rows(4).Insert()

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 25236.75
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 2025
$ws.Range("E4").Value = "06/2025"
